# Re-shuffle the per-row market data (Fecha / Volumen / Precio minimo /
# Precio maximo / Precio promedio ponderado / Precio $/Kg) across the
# existing data rows (2-41), turning the "daily" logic into the
# "weekly" logic. Only columns D, J, K, L, M, P move; everything else
# (market id, region, category, quality, unit, origin, etc.) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# final row -> source row (which row's D/J/K/L/M/P values end up here)
$mapping = @{
    2  = 23
    3  = 38
    4  = 3
    5  = 5
    6  = 8
    7  = 39
    8  = 20
    9  = 26
    10 = 7
    11 = 41
    12 = 21
    13 = 24
    14 = 19
    15 = 31
    16 = 27
    17 = 2
    18 = 6
    19 = 29
    20 = 36
    21 = 11
    22 = 33
    23 = 30
    24 = 40
    25 = 4
    26 = 17
    27 = 16
    28 = 10
    29 = 18
    30 = 35
    31 = 28
    32 = 32
    33 = 9
    34 = 14
    35 = 13
    36 = 22
    37 = 15
    38 = 37
    39 = 25
    40 = 12
    41 = 34
}

# Columns that move together as a group (by index: D=4, J=10, K=11, L=12, M=13, P=16)
$cols = @(4, 10, 11, 12, 13, 16)

# 1) Snapshot the original values for every data row before writing anything,
#    so the permutation (which includes multi-element cycles, not just swaps)
#    doesn't clobber a source row before it has been read.
$snapshot = @{}
for ($row = 2; $row -le 41; $row++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowVals
}

# 2) Write each row's new values from the snapshot of its mapped source row.
for ($row = 2; $row -le 41; $row++) {
    $srcRow = $mapping[$row]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($row, $col).Value = $srcVals[$col]
    }
}
